# CIERRE 4 DE DIC 2021
# Advance the payroll receipt from "SEMANA 48 (22-28 NOV 2021)" to
# "SEMANA 49 (29 NOV - 05 DIC 2021)" and update the corresponding figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# --- Week banner (drives every =B9 / =B27 / =H27 / =B43 / =H43 / =B60 formula) ---
$ws.Range("B9").Value2 = "SEMANA   49  DEL    29      Al   05   DE   DICIEMBRE          2021"

# --- First employee block (rows 3-6) ---
$ws.Range("J3").Value2 = 5
$ws.Range("K3").Value2 = 2167
$ws.Range("K4").Value2 = 0

# --- Second employee block (rows 38-41) ---
$ws.Range("D38").Value2 = 5
$ws.Range("E38").Value2 = 1833
$ws.Range("K39").Value2 = 833

# --- View state: scrolled down to show the second receipt, C58 selected ---
$ws.Range("C58").Select()
